# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the inlineStr cells already used
# throughout this sheet). Without forcing the format, Excel auto-detects
# numeric-looking strings (e.g. "0.999", "538.96") and stores them as
# numbers instead of text, which would not match the source data.
function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Price (D) / Volume 1h (E) refresh for existing coins ---
Set-TextValue $ws 'D2' '59.433.01'
Set-TextValue $ws 'E2' '  +0.65%  '
Set-TextValue $ws 'D3' '2.604.15'
Set-TextValue $ws 'E3' '  +0.66%  '
Set-TextValue $ws 'E4' '  -0.79%  '
Set-TextValue $ws 'D5' '538.96'
Set-TextValue $ws 'E5' '  +3.17%  '
Set-TextValue $ws 'D6' '141.52'
Set-TextValue $ws 'E6' '  +1.57%  '
Set-TextValue $ws 'D7' '0.999'
Set-TextValue $ws 'E7' '  +0.08%  '
Set-TextValue $ws 'D8' '0.565'
Set-TextValue $ws 'E8' '  +0.04%  '
Set-TextValue $ws 'D9' '6.51'
Set-TextValue $ws 'E9' '  -0.39%  '
Set-TextValue $ws 'E10' '  +1.10%  '
Set-TextValue $ws 'E11' '  +1.39%  '
Set-TextValue $ws 'E12' '  -1.02%  '
Set-TextValue $ws 'D13' '3.060.75'
Set-TextValue $ws 'E13' '  +0.35%  '
Set-TextValue $ws 'D14' '59.370.13'
Set-TextValue $ws 'E14' '  +0.67%  '
Set-TextValue $ws 'D15' '20.77'
Set-TextValue $ws 'E15' '  +1.11%  '
Set-TextValue $ws 'D16' '2.642.45'
Set-TextValue $ws 'E16' '  +1.52%  '
Set-TextValue $ws 'E17' '  +0.35%  '
Set-TextValue $ws 'D18' '340.99'
Set-TextValue $ws 'E18' '  -0.04%  '
Set-TextValue $ws 'E19' '  +1.35%  '
Set-TextValue $ws 'D20' '10.09'
Set-TextValue $ws 'E20' '  +0.04%  '
Set-TextValue $ws 'D21' '6.31'
Set-TextValue $ws 'E21' '  -2.19%  '
Set-TextValue $ws 'E22' '  -0.02%  '
Set-TextValue $ws 'D23' '67.19'
Set-TextValue $ws 'E23' '  +0.97%  '
Set-TextValue $ws 'D24' '0.409'
Set-TextValue $ws 'E24' '  +1.29%  '
Set-TextValue $ws 'E25' '  -1.29%  '
Set-TextValue $ws 'D26' '0.999'
Set-TextValue $ws 'E26' '  +0.08%  '
Set-TextValue $ws 'D27' '7.22'
Set-TextValue $ws 'E27' '  +2.08%  '
Set-TextValue $ws 'D28' '0.0₃0744'
Set-TextValue $ws 'E28' '  +2.68%  '
Set-TextValue $ws 'E29' '  +0.04%  '
Set-TextValue $ws 'E30' '  +6.12%  '
Set-TextValue $ws 'D31' '5.84'
Set-TextValue $ws 'E31' '  -1.02%  '
Set-TextValue $ws 'D32' '18.79'
Set-TextValue $ws 'E32' '  +0.32%  '
Set-TextValue $ws 'D33' '150.32'
Set-TextValue $ws 'E33' '  +0.73%  '
Set-TextValue $ws 'D34' '3.99'
Set-TextValue $ws 'E34' '  +0.63%  '
Set-TextValue $ws 'E35' '  +0.20%  '
Set-TextValue $ws 'D36' '0.847'
Set-TextValue $ws 'E36' '  +4.53%  '
Set-TextValue $ws 'E37' '  -0.63%  '
Set-TextValue $ws 'D38' '0.825'
Set-TextValue $ws 'E38' '  -0.07%  '
Set-TextValue $ws 'E39' '  +0.17%  '
Set-TextValue $ws 'D40' '0.999'
Set-TextValue $ws 'E40' '  +0.16%  '
Set-TextValue $ws 'D41' '274.85'
Set-TextValue $ws 'E41' '  +0.81%  '
Set-TextValue $ws 'E42' '  -0.69%  '
Set-TextValue $ws 'D43' '10.72'
Set-TextValue $ws 'E43' '  -0.21%  '
Set-TextValue $ws 'E44' '  -0.08%  '
Set-TextValue $ws 'E45' '  +1.25%  '
Set-TextValue $ws 'D48' '1.939.01'
Set-TextValue $ws 'E48' '  -1.53%  '
Set-TextValue $ws 'E49' '  -0.02%  '
Set-TextValue $ws 'D50' '111.46'
Set-TextValue $ws 'E50' '  -2.00%  '
Set-TextValue $ws 'E51' '  +2.12%  '

# --- Rows 46/47: coin ranking swapped (VeChain now ranks above InjectiveProtocol) ---
Set-TextValue $ws 'B46' 'VeChain'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D46' '0.0223'
Set-TextValue $ws 'E46' '  +1.21%  '

Set-TextValue $ws 'B47' 'InjectiveProtocol'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D47' '18.49'
Set-TextValue $ws 'E47' '  +2.91%  '

Write-Host "Applied cryptos.xlsx update"
